# Auto-generated edit script: update Price (D) and Volume(1h) (E) columns
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @("25.862.29", "  -1.06%  ")
    3 = @("1.635.95", "  -0.89%  ")
    4 = @("1.002", "  -0.32%  ")
    5 = @("214.81", "  -0.17%  ")
    6 = @("0.5019", "  -1.77%  ")
    7 = @($null, "  -0.34%  ")
    8 = @("0.2561", "  -1.07%  ")
    9 = @("0.06377", "  -0.86%  ")
    10 = @("19.68", "  -1.20%  ")
    11 = @("0.07706", "  -1.08%  ")
    12 = @("1.657.64", "  +0.30%  ")
    13 = @("4.261", "  -0.47%  ")
    14 = @("1.860.66", "  -0.93%  ")
    15 = @("0.5442", "  -1.32%  ")
    16 = @("0.0₅7899", "  -1.32%  ")
    17 = @("64.17", $null)
    18 = @("25.854.94", "  -1.13%  ")
    19 = @($null, "  -0.24%  ")
    20 = @("203.25", "  -3.62%  ")
    21 = @("4.375", "  -0.27%  ")
    22 = @("9.895", "  -1.61%  ")
    23 = @("5.970", "  -1.08%  ")
    24 = @($null, "  -0.28%  ")
    25 = @("1.919", "  +9.58%  ")
    26 = @("141.21", $null)
    27 = @("0.1136", "  -3.31%  ")
    28 = @("15.65", "  -0.93%  ")
    29 = @("6.712", "  -3.79%  ")
    30 = @("1.241", "  +0.00%  ")
    31 = @($null, "  -4.00%  ")
    32 = @("3.276", "  -2.12%  ")
    33 = @("3.187", "  -0.77%  ")
    34 = @("1.540", "  -1.11%  ")
    35 = @("2.370", "  +0.94%  ")
    36 = @("2.625", "  -4.16%  ")
    37 = @("0.8918", "  -3.40%  ")
    38 = @("1.158.75", "  -0.88%  ")
    39 = @("0.5603", "  -1.52%  ")
    40 = @("0.01561", "  -1.47%  ")
    41 = @("1.001", "  -0.35%  ")
    42 = @("5.706", "  +0.90%  ")
    43 = @("0.8075", "  -1.85%  ")
    44 = @("99.62", "  -0.32%  ")
    45 = @("1.773.44", $null)
    46 = @($null, "  -0.74%  ")
    47 = @("0.4514", "  -0.82%  ")
    48 = @("1.003", "  -0.25%  ")
    49 = @("54.97", "  -0.86%  ")
    50 = @("0.05053", $null)
    51 = @("1.002", "  -0.46%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals[0]) {
        $ws.Cells.Item([int]$row, 4).Value = $vals[0]
    }
    if ($null -ne $vals[1]) {
        $ws.Cells.Item([int]$row, 5).Value = $vals[1]
    }
}

